# Correcao no regex/mensagem da data: troca a mensagem de formato de data
# de "YYYY-MM-DD" para "DD-MM-AAAA" nas linhas que validam campos de data
# (data_nascimento, data_preenchimento, data_validacao, data_recebimento).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "Data no formato DD-MM-AAAA (ex: 30-05-2024)"

$ws.Range("D10").Value = $newText
$ws.Range("D63").Value = $newText
$ws.Range("D64").Value = $newText
$ws.Range("D65").Value = $newText

# Deixa a selecao na ultima celula editada, como ficou no ficheiro de origem
$ws.Range("D65").Select()
